$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 9500
$ws.Range("J26").Value = 9500
$ws.Range("L26").Value = 9500
$ws.Range("N26").Value = -10188

$ws.Range("H62").Value = 2534.1667
$ws.Range("I62").Value = 2489.375
$ws.Range("J62").Value = 2623.75
$ws.Range("K62").Value = 2489.375
$ws.Range("L62").Value = 2623.75
$ws.Range("M62").Value = -1865.375
$ws.Range("N62").Value = -3871.75

$ws.Range("H65").Value = 2534.1667
$ws.Range("I65").Value = 2489.375
$ws.Range("J65").Value = 2623.75
$ws.Range("K65").Value = 12446.875
$ws.Range("L65").Value = 13118.75
$ws.Range("M65").Value = -9326.875
$ws.Range("N65").Value = -19358.75

$ws.Range("H98").Value = 1134.9445
$ws.Range("I98").Value = 1194
$ws.Range("J98").Value = 981.4
$ws.Range("K98").Value = 1194
$ws.Range("L98").Value = 981.4
$ws.Range("M98").Value = 304
$ws.Range("N98").Value = -3977.4

$ws.Range("H100").Value = 1752.5
$ws.Range("I100").Value = 1005
$ws.Range("K100").Value = 1005
$ws.Range("M100").Value = -464

$ws.Range("H122").Value = 1134.9445
$ws.Range("I122").Value = 1194
$ws.Range("J122").Value = 981.4
$ws.Range("K122").Value = 3582
$ws.Range("L122").Value = 2944.2
$ws.Range("M122").Value = -1132
$ws.Range("N122").Value = -7844.2

$ws.Range("H129").Value = 866.6721
$ws.Range("I129").Value = 319.4
$ws.Range("J129").Value = 915.5357
$ws.Range("K129").Value = 958.1999999999999
$ws.Range("L129").Value = 2746.6071
$ws.Range("M129").Value = 4041.8
$ws.Range("N129").Value = -12746.6071

$ws.Range("H138").Value = 2943.4607
$ws.Range("I138").Value = 1514.6285
$ws.Range("J138").Value = 3869.5557
$ws.Range("K138").Value = 4543.8855
$ws.Range("L138").Value = 11608.6671
$ws.Range("M138").Value = 596.1144999999997
$ws.Range("N138").Value = -21888.6671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H119").Value = 60698
$ws.Range("J119").Value = 60698
$ws.Range("L119").Value = 60698
$ws.Range("N119").Value = -70374

$ws.Range("H121").Value = 46200
$ws.Range("J121").Value = 46200
$ws.Range("L121").Value = 46200
$ws.Range("N121").Value = -49694

$ws.Range("H122").Value = 1603.6875
$ws.Range("I122").Value = 1353.4
$ws.Range("J122").Value = 2020.8334
$ws.Range("K122").Value = 4060.2
$ws.Range("L122").Value = 6062.5002
$ws.Range("M122").Value = -1610.2
$ws.Range("N122").Value = -10962.5002

$ws.Range("H132").Value = 5125.1
$ws.Range("I132").Value = 2815.7693
$ws.Range("J132").Value = 7626.875
$ws.Range("K132").Value = 8447.3079
$ws.Range("L132").Value = 22880.625
$ws.Range("M132").Value = -5917.3079
$ws.Range("N132").Value = -27940.625

$ws.Range("H134").Value = 61079.855
$ws.Range("J134").Value = 61079.855
$ws.Range("L134").Value = 61079.855
$ws.Range("N134").Value = -71219.85500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2869.3333
$ws.Range("I20").Value = 1608
$ws.Range("J20").Value = 3500
$ws.Range("K20").Value = 1608
$ws.Range("L20").Value = 3500
$ws.Range("M20").Value = -1361
$ws.Range("N20").Value = -3994

$ws.Range("H132").Value = 59506.668
$ws.Range("J132").Value = 59506.668
$ws.Range("L132").Value = 59506.668
$ws.Range("N132").Value = -69626.66800000001

$ws.Range("H134").Value = 3280.2666
$ws.Range("I134").Value = 2683.6667
$ws.Range("J134").Value = 5666.6665
$ws.Range("K134").Value = 8051.000100000001
$ws.Range("L134").Value = 16999.9995
$ws.Range("M134").Value = -5516.000100000001
$ws.Range("N134").Value = -22069.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 12144.143
$ws.Range("I32").Value = 4477.25
$ws.Range("J32").Value = 22366.666
$ws.Range("K32").Value = 4477.25
$ws.Range("L32").Value = 22366.666
$ws.Range("M32").Value = -4161.25
$ws.Range("N32").Value = -22998.666

$ws.Range("H99").Value = 2319.7144
$ws.Range("I99").Value = 1952.4
$ws.Range("J99").Value = 2523.7778
$ws.Range("K99").Value = 1952.4
$ws.Range("L99").Value = 2523.7778
$ws.Range("M99").Value = -454.4000000000001
$ws.Range("N99").Value = -5519.7778

$ws.Range("H126").Value = 2319.7144
$ws.Range("I126").Value = 1952.4
$ws.Range("J126").Value = 2523.7778
$ws.Range("K126").Value = 5857.200000000001
$ws.Range("L126").Value = 7571.3334
$ws.Range("M126").Value = -3387.200000000001
$ws.Range("N126").Value = -12511.3334

$ws.Range("H134").Value = 3888.1035
$ws.Range("I134").Value = 3716.1
$ws.Range("J134").Value = 3978.6316
$ws.Range("K134").Value = 11148.3
$ws.Range("L134").Value = 11935.8948
$ws.Range("M134").Value = -8613.299999999999
$ws.Range("N134").Value = -17005.8948

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 5925
$ws.Range("I130").Value = 1250
$ws.Range("K130").Value = 3750
$ws.Range("M130").Value = 1270

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9337.5
$ws.Range("I80").Value = 26000
$ws.Range("J80").Value = 3783.3333
$ws.Range("K80").Value = 26000
$ws.Range("L80").Value = 3783.3333
$ws.Range("M80").Value = -25002
$ws.Range("N80").Value = -5779.3333

$ws.Range("H83").Value = 9337.5
$ws.Range("I83").Value = 26000
$ws.Range("J83").Value = 3783.3333
$ws.Range("K83").Value = 130000
$ws.Range("L83").Value = 18916.6665
$ws.Range("M83").Value = -125008
$ws.Range("N83").Value = -28900.6665

$ws.Range("H102").Value = 4060.2222
$ws.Range("I102").Value = 2929.0435
$ws.Range("J102").Value = 6061.5386
$ws.Range("K102").Value = 2929.0435
$ws.Range("L102").Value = 6061.5386
$ws.Range("M102").Value = -1307.0435
$ws.Range("N102").Value = -9305.5386

$ws.Range("H122").Value = 2527
$ws.Range("I122").Value = 2834.1428
$ws.Range("J122").Value = 1571.4445
$ws.Range("K122").Value = 8502.428400000001
$ws.Range("L122").Value = 4714.333500000001
$ws.Range("M122").Value = -6052.428400000001
$ws.Range("N122").Value = -9614.333500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3326.9688
$ws.Range("I40").Value = 3106.7917
$ws.Range("J40").Value = 3987.5
$ws.Range("K40").Value = 3106.7917
$ws.Range("L40").Value = 3987.5
$ws.Range("M40").Value = -2970.7917
$ws.Range("N40").Value = -4259.5

$ws.Range("H82").Value = 2946
$ws.Range("I82").Value = 3500
$ws.Range("J82").Value = 2761.3333
$ws.Range("K82").Value = 3500
$ws.Range("L82").Value = 2761.3333
$ws.Range("M82").Value = -3139
$ws.Range("N82").Value = -3483.3333

$ws.Range("H85").Value = 2946
$ws.Range("I85").Value = 3500
$ws.Range("J85").Value = 2761.3333
$ws.Range("K85").Value = 3500
$ws.Range("L85").Value = 2761.3333
$ws.Range("M85").Value = -2252
$ws.Range("N85").Value = -5257.3333

$ws.Range("H122").Value = 4660.2925
$ws.Range("I122").Value = 4468.4
$ws.Range("J122").Value = 5299.933
$ws.Range("K122").Value = 13405.2
$ws.Range("L122").Value = 15899.799
$ws.Range("M122").Value = -10955.2
$ws.Range("N122").Value = -20799.799

$ws.Range("H132").Value = 3336.034
$ws.Range("I132").Value = 2960.6667
$ws.Range("J132").Value = 4263.4116
$ws.Range("K132").Value = 8882.000100000001
$ws.Range("L132").Value = 12790.2348
$ws.Range("M132").Value = -6352.000100000001
$ws.Range("N132").Value = -17850.2348

$ws.Range("H136").Value = 5539.326
$ws.Range("I136").Value = 4228.56
$ws.Range("J136").Value = 7099.7617
$ws.Range("K136").Value = 12685.68
$ws.Range("L136").Value = 21299.2851
$ws.Range("M136").Value = -10135.68
$ws.Range("N136").Value = -26399.2851

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H140").Value = 58046.273
$ws.Range("J140").Value = 58046.273
$ws.Range("L140").Value = 58046.273
$ws.Range("N140").Value = -68406.273
